$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.58
$ws.Range("P2").Value = 1.96
$ws.Range("G3").Value = 3.75
$ws.Range("I3").Value = 2.56
$ws.Range("J3").Value = 3.5
$ws.Range("G4").Value = 1.44
$ws.Range("I4").Value = 12
$ws.Range("Q4").Value = 1.56
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.2
$ws.Range("G5").Value = 3.4
$ws.Range("J5").Value = 3.5
$ws.Range("F6").Value = 2.26
$ws.Range("G6").Value = 2.56
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3.5
$ws.Range("P6").Value = 1.62
$ws.Range("Q6").Value = 2.28
$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 9.199999999999999
$ws.Range("L7").Value = 1.47
$ws.Range("N7").Value = 2.92
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 1.64
$ws.Range("Q7").Value = 2.16
$ws.Range("T7").Value = 2.66
$ws.Range("U7").Value = 1.52
$ws.Range("X7").Value = 12.5
$ws.Range("AB7").Value = 5.7
$ws.Range("AN7").Value = 980
$ws.Range("Q8").Value = 2.94
$ws.Range("P9").Value = 1.24
$ws.Range("F10").Value = 1.54
$ws.Range("G10").Value = 1.76
$ws.Range("H10").Value = 2.4
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 3.65
$ws.Range("P10").Value = 1.68
$ws.Range("Q10").Value = 1.99
$ws.Range("G11").Value = 2.44
$ws.Range("I11").Value = 4.4
$ws.Range("J11").Value = 3
$ws.Range("P11").Value = 1.53
$ws.Range("Q11").Value = 2.68
$ws.Range("G12").Value = 1.54
$ws.Range("F13").Value = 1.39
$ws.Range("G13").Value = 1.48
$ws.Range("H13").Value = 9
$ws.Range("J13").Value = 4.6
$ws.Range("P13").Value = 2.02
$ws.Range("Q13").Value = 1.8
